$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the number-format style (s="1" -> default) from A1, B1, A2, B2 first,
# so the subsequent value write isn't coerced to text by the old Text format
# (numFmtId 49) associated with style index 1.
$ws.Range("A1:B2").Style = "Normal"

# Update B1's value (708 -> 190)
$ws.Range("B1").Value = 190

# Update the active selection to just C4 (was A1:C4)
$ws.Range("C4").Select()
